# Weekly update: insert two new "Haba" price records for Femacal de La Calera
# right after the existing row 146, pushing the remaining historical rows
# (old 147-156) down by two rows (new 149-158).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 147; existing rows 147:156 shift to 149:158.
$ws.Rows("147:148").Insert()

# New row 147
$ws.Cells.Item(147, 1).Value = 3
$ws.Cells.Item(147, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(147, 3).Value = "Coquimbo"
$ws.Cells.Item(147, 4).Value = 44753
$ws.Cells.Item(147, 5).Value = 5
$ws.Cells.Item(147, 6).Value = 100112026
$ws.Cells.Item(147, 7).Value = "Haba"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 80
$ws.Cells.Item(147, 11).Value = 17000
$ws.Cells.Item(147, 12).Value = 18000
$ws.Cells.Item(147, 13).Value = 17500
$ws.Cells.Item(147, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(147, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(147, 16).Value = 700
$ws.Cells.Item(147, 17).Value = 25
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# New row 148
$ws.Cells.Item(148, 1).Value = 3
$ws.Cells.Item(148, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(148, 3).Value = "Coquimbo"
$ws.Cells.Item(148, 4).Value = 44753
$ws.Cells.Item(148, 5).Value = 5
$ws.Cells.Item(148, 6).Value = 100112026
$ws.Cells.Item(148, 7).Value = "Haba"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Segunda"
$ws.Cells.Item(148, 10).Value = 45
$ws.Cells.Item(148, 11).Value = 14000
$ws.Cells.Item(148, 12).Value = 14000
$ws.Cells.Item(148, 13).Value = 14000
$ws.Cells.Item(148, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(148, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(148, 16).Value = 560
$ws.Cells.Item(148, 17).Value = 25
$ws.Cells.Item(148, 18).Value = "Hortaliza"
